# Update the NCAP_BND assumptions in rows 12-21 (column E).
# The original formulas scaled each figure by 3/5 (e.g. "=13*3/5"); the
# revised workbook drops that scaling and uses the raw figures instead
# (some cells keep a trivial "=N" formula, two become plain literal values).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NCAP_BND")

$ws.Range("E12").Formula = "=13"
$ws.Range("E13").Formula = "=18"
$ws.Range("E14").Formula = "=25"
$ws.Range("E15").Value = 34
$ws.Range("E16").Formula = "=40"

$ws.Range("E17").Formula = "=8.4"
$ws.Range("E18").Formula = "=12"
$ws.Range("E19").Formula = "=18"
$ws.Range("E20").Formula = "=25"
$ws.Range("E21").Value = 33

# Mirror the saved cursor/selection position left behind by the edit.
$ws.Range("F12").Select()
